$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "72.654.50"
$ws.Range("E2").Value = "  +5.56%  "
$ws.Range("D3").Value = "4.063.72"
$ws.Range("E3").Value = "  +5.63%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "'521.25"
$ws.Range("E5").Value = "  -0.41%  "
$ws.Range("D6").Value = "'147.78"
$ws.Range("E6").Value = "  +4.13%  "
$ws.Range("D7").Value = "'0.729"
$ws.Range("E7").Value = "  +20.45%  "
$ws.Range("D8").Value = "4.053.87"
$ws.Range("E8").Value = "  +5.69%  "
$ws.Range("E9").Value = "  +0.16%  "
$ws.Range("D10").Value = "'0.783"
$ws.Range("E10").Value = "  +10.19%  "
$ws.Range("D11").Value = "'0.176"
$ws.Range("E11").Value = "  +4.51%  "
$ws.Range("D12").Value = "'0.0000331"
$ws.Range("E12").Value = "  +0.98%  "
$ws.Range("D13").Value = "'48.70"
$ws.Range("E13").Value = "  +17.32%  "
$ws.Range("D14").Value = "'11.12"
$ws.Range("E14").Value = "  +10.10%  "
$ws.Range("D15").Value = "4.707.57"
$ws.Range("E15").Value = "  +5.48%  "
$ws.Range("D16").Value = "4.058.56"
$ws.Range("E16").Value = "  +5.18%  "
$ws.Range("B17").Value = "Uniswap"
$ws.Range("C17").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D17").Value = "'14.41"
$ws.Range("E17").Value = "  +4.31%  "
$ws.Range("B18").Value = "Chainlink"
$ws.Range("C18").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D18").Value = "'21.34"
$ws.Range("E18").Value = "  +3.36%  "
$ws.Range("E19").Value = "  +2.42%  "
$ws.Range("E20").Value = "  -0.07%  "
$ws.Range("D21").Value = "72.531.08"
$ws.Range("E21").Value = "  +5.60%  "
$ws.Range("D22").Value = "'448.41"
$ws.Range("E22").Value = "  +7.08%  "
$ws.Range("D23").Value = "'105.24"
$ws.Range("E23").Value = "  +21.38%  "
$ws.Range("D24").Value = "'3.60"
$ws.Range("E24").Value = "  +6.95%  "
$ws.Range("D25").Value = "'15.11"
$ws.Range("E25").Value = "  +7.88%  "
$ws.Range("E26").Value = "  +2.18%  "
$ws.Range("E27").Value = "  +1.27%  "
$ws.Range("D28").Value = "'11.07"
$ws.Range("E28").Value = "  +5.70%  "
$ws.Range("D29").Value = "'38.12"
$ws.Range("E29").Value = "  +6.31%  "
$ws.Range("D30").Value = "'5.83"
$ws.Range("E30").Value = "  +2.83%  "
$ws.Range("D31").Value = "'3.28"
$ws.Range("E31").Value = "  +16.26%  "
$ws.Range("D32").Value = "'13.69"
$ws.Range("E32").Value = "  +4.99%  "
$ws.Range("D33").Value = "'0.130"
$ws.Range("E33").Value = "  +4.33%  "
$ws.Range("D34").Value = "'677.26"
$ws.Range("E34").Value = "  -0.95%  "
$ws.Range("D35").Value = "'68.15"
$ws.Range("E35").Value = "  +0.89%  "
$ws.Range("E36").Value = "  +12.97%  "
$ws.Range("D37").Value = "'42.20"
$ws.Range("E37").Value = "  +7.01%  "
$ws.Range("B38").Value = "TheGraph"
$ws.Range("C38").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D38").Value = "'0.431"
$ws.Range("E38").Value = "  +0.09%  "
$ws.Range("B39").Value = "PEPE"
$ws.Range("C39").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D39").Value = "0.0₃0862"
$ws.Range("E39").Value = "  +2.45%  "
$ws.Range("E40").Value = "  +4.40%  "
$ws.Range("E41").Value = "  +8.58%  "
$ws.Range("D43").Value = "'0.0499"
$ws.Range("E43").Value = "  +4.97%  "
$ws.Range("D44").Value = "'0.999"
$ws.Range("E44").Value = "  -0.27%  "
$ws.Range("E45").Value = "  +2.50%  "
$ws.Range("E46").Value = "  +14.06%  "
$ws.Range("B47").Value = "THORChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D47").Value = "'9.85"
$ws.Range("E47").Value = "  +16.81%  "
$ws.Range("B48").Value = "Fetch.AI"
$ws.Range("C48").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D48").Value = "'2.68"
$ws.Range("E48").Value = "  -1.59%  "
$ws.Range("D49").Value = "'3.44"
$ws.Range("E49").Value = "  +1.72%  "
$ws.Range("D50").Value = "'3.08"
$ws.Range("E50").Value = "  +5.17%  "
